$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, D(date serial), H(variedad), J(volumen), K(precio min), L(precio max), M(precio promedio), P(precio $/Kg)
$data = @(
    @(2, 44741, "Perfection", 160, 28000, 30000, 29000, 1160),
    @(3, 44818, "Perfection", 400, 24000, 27000, 25500, 1020),
    @(4, 45126, "Perfection", 600, 22000, 24000, 23000, 920),
    @(5, 44783, "Perfection", 400, 27000, 29000, 28000, 1120),
    @(6, 44447, "Perfection", 600, 28000, 30000, 29000, 1160),
    @(7, 45141, "Perfection", 240, 19000, 20000, 19500, 780),
    @(8, 44812, "Perfection", 500, 28000, 30000, 29000, 1160),
    @(9, 44357, "Perfection", 340, 28000, 30000, 29000, 1160),
    @(10, 45106, "Perfection", 900, 26000, 28000, 27000, 1080),
    @(11, 44343, "Perfection", 200, 26000, 28000, 27000, 1080),
    @(12, 44448, "Perfection", 400, 28000, 30000, 29000, 1160),
    @(13, 44755, "Perfection", 200, 30000, 32000, 31000, 1240),
    @(14, 44671, "Sin especificar", 240, 23000, 25000, 24000, 960),
    @(15, 44482, "Perfection", 500, 18000, 20000, 19000, 760),
    @(16, 44720, "Perfection", 400, 28000, 30000, 29000, 1160),
    @(17, 44798, "Perfection", 400, 30000, 32000, 31000, 1240),
    @(18, 44497, "Perfection", 500, 13000, 15000, 14000, 560),
    @(19, 45120, "Perfection", 1100, 25000, 27000, 26000, 1040),
    @(20, 44874, "Perfection", 160, 14000, 16000, 15000, 600),
    @(21, 44791, "Perfection", 500, 27000, 29000, 28000, 1120),
    @(22, 45155, "Perfection", 1000, 23000, 25000, 24000, 960),
    @(23, 44468, "Perfection", 500, 23000, 25000, 24000, 960),
    @(24, 45148, "Perfection", 600, 20000, 22000, 21000, 840),
    @(25, 44469, "Perfection", 600, 22000, 24000, 23000, 920),
    @(26, 44406, "Perfection", 600, 26000, 28000, 27000, 1080),
    @(27, 44356, "Perfection", 300, 26000, 28000, 27000, 1080),
    @(28, 45084, "Perfection", 500, 25000, 27000, 26000, 1040),
    @(29, 44427, "Perfection", 300, 28000, 30000, 29000, 1160),
    @(30, 44349, "Perfection", 600, 26000, 28000, 27000, 1080),
    @(31, 44391, "Perfection", 100, 26000, 28000, 27000, 1080),
    @(32, 44742, "Perfection", 200, 28000, 30000, 29000, 1160),
    @(33, 44825, "Perfection", 480, 28000, 30000, 29000, 1160),
    @(34, 44398, "Perfection", 500, 26000, 28000, 27000, 1080),
    @(35, 45161, "Perfection", 800, 22000, 24000, 23000, 920),
    @(36, 45204, "Sin especificar", 400, 20000, 21000, 20500, 820),
    @(37, 45113, "Perfection", 800, 23000, 25000, 24000, 960),
    @(38, 45175, "Perfection", 400, 21000, 23000, 22000, 880),
    @(39, 44490, "Perfection", 500, 16000, 18000, 17000, 680),
    @(40, 45169, "Perfection", 700, 21000, 23000, 22000, 880),
    @(41, 44489, "Perfection", 400, 18000, 20000, 19000, 760),
    @(42, 45134, "Perfection", 600, 23000, 25000, 24000, 960),
    @(43, 44685, "Perfection", 160, 25000, 27000, 26000, 1040),
    @(44, 44748, "Perfection", 700, 28000, 30000, 29000, 1160),
    @(45, 44384, "Perfection", 400, 26000, 28000, 27000, 1080),
    @(46, 44706, "Perfection", 160, 25000, 26000, 25500, 1020),
    @(47, 45133, "Perfection", 560, 23000, 25000, 24000, 960),
    @(48, 44839, "Perfection", 700, 22000, 24000, 23000, 920),
    @(49, 44364, "Perfection", 200, 28000, 30000, 29000, 1160),
    @(50, 44707, "Perfection", 200, 25000, 27000, 26000, 1040),
    @(51, 44678, "Perfection", 600, 25000, 27000, 26000, 1040),
    @(52, 44435, "Perfection", 900, 28000, 30000, 29000, 1160),
    @(53, 44679, "Perfection", 400, 25000, 27000, 26000, 1040),
    @(54, 45218, "Perfection", 400, 20000, 22000, 21000, 840),
    @(55, 44811, "Perfection", 700, 28000, 30000, 29000, 1160),
    @(56, 45092, "Perfection", 300, 27000, 29000, 28000, 1120),
    @(57, 45119, "Perfection", 1000, 26000, 28000, 27000, 1080),
    @(58, 44714, "Perfection", 240, 27000, 28000, 27500, 1100),
    @(59, 44868, "Perfection", 300, 13000, 15000, 14000, 560),
    @(60, 45112, "Perfection", 1000, 25000, 26000, 25500, 1020),
    @(61, 44420, "Perfection", 700, 27000, 29000, 28000, 1120),
    @(62, 44441, "Perfection", 700, 28000, 30000, 29000, 1160),
    @(63, 44756, "Perfection", 240, 30000, 32000, 31000, 1240),
    @(64, 45127, "Perfection", 700, 20000, 22000, 21000, 840),
    @(65, 44426, "Perfection", 400, 28000, 30000, 29000, 1160),
    @(66, 44377, "Perfection", 500, 26000, 28000, 27000, 1080),
    @(67, 44405, "Perfection", 500, 26000, 28000, 27000, 1080),
    @(68, 44784, "Perfection", 360, 27000, 29000, 28000, 1120),
    @(69, 44763, "Perfection", 400, 29000, 30000, 29500, 1180),
    @(70, 44363, "Perfection", 240, 28000, 30000, 29000, 1160),
    @(71, 44370, "Perfection", 400, 27000, 28000, 27500, 1100),
    @(72, 44475, "Perfection", 1000, 22000, 24000, 23000, 920),
    @(73, 44462, "Perfection", 400, 22000, 23000, 22500, 900),
    @(74, 44769, "Perfection", 500, 30000, 32000, 31000, 1240),
    @(75, 44350, "Perfection", 700, 28000, 30000, 29000, 1160),
    @(76, 44454, "Perfection", 1000, 28000, 30000, 29000, 1160),
    @(77, 44483, "Perfection", 300, 18000, 20000, 19000, 760),
    @(78, 45196, "Perfection", 400, 19000, 20000, 19500, 780),
    @(79, 44413, "Perfection", 700, 26000, 28000, 27000, 1080),
    @(80, 44790, "Perfection", 560, 27000, 29000, 28000, 1120),
    @(81, 44721, "Perfection", 240, 28000, 30000, 29000, 1160),
    @(82, 44399, "Perfection", 400, 26000, 28000, 27000, 1080),
    @(83, 45085, "Perfection", 400, 25000, 27000, 26000, 1040),
    @(84, 44846, "Sin especificar", 488, 23000, 24000, 23426, 937),
    @(85, 45140, "Perfection", 300, 20000, 22000, 21000, 840),
    @(86, 44727, "Perfection", 160, 28000, 30000, 29000, 1160),
    @(87, 44476, "Perfection", 500, 23000, 24000, 23500, 940),
    @(88, 44762, "Perfection", 400, 29000, 30000, 29500, 1180),
    @(89, 45217, "Perfection", 600, 21000, 23000, 22000, 880),
    @(90, 45162, "Perfection", 700, 22000, 23000, 22500, 900),
    @(91, 44412, "Perfection", 600, 25000, 27000, 26000, 1040),
    @(92, 45176, "Perfection", 240, 21000, 23000, 22000, 880),
    @(93, 45147, "Perfection", 700, 20000, 22000, 21000, 840),
    @(94, 44434, "Perfection", 500, 28000, 30000, 29000, 1160),
    @(95, 44672, "Sin especificar", 160, 23000, 25000, 24000, 960),
    @(96, 44371, "Perfection", 500, 28000, 30000, 29000, 1160),
    @(97, 44392, "Perfection", 100, 26000, 28000, 27000, 1080),
    @(98, 44819, "Perfection", 500, 25000, 28000, 26500, 1060),
    @(99, 44749, "Perfection", 470, 28000, 30000, 29064, 1163),
    @(101, 44433, "Perfection", 400, 28000, 30000, 29000, 1160),
    @(102, 44847, "Sin especificar", 300, 23000, 24000, 23500, 940),
    @(103, 44419, "Perfection", 600, 27000, 29000, 28000, 1120),
    @(104, 45168, "Perfection", 600, 21000, 23000, 22000, 880),
    @(105, 44804, "Perfection", 400, 28000, 30000, 29000, 1160),
    @(106, 44776, "Perfection", 400, 28000, 30000, 29000, 1160),
    @(107, 44860, "Perfection", 200, 15000, 16000, 15500, 620),
    @(108, 44455, "Perfection", 800, 28000, 30000, 29000, 1160),
    @(109, 44826, "Perfection", 520, 28000, 30000, 29000, 1160),
    @(110, 44461, "Perfection", 500, 23000, 25000, 24000, 960),
    @(111, 44385, "Perfection", 500, 26000, 28000, 27000, 1080),
    @(112, 45091, "Perfection", 360, 26000, 28000, 27000, 1080),
    @(113, 44699, "Perfection", 200, 29000, 30000, 29500, 1180)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 8).Value = $row[2]
    $ws.Cells.Item($r, 10).Value = $row[3]
    $ws.Cells.Item($r, 11).Value = $row[4]
    $ws.Cells.Item($r, 12).Value = $row[5]
    $ws.Cells.Item($r, 13).Value = $row[6]
    $ws.Cells.Item($r, 16).Value = $row[7]
}

Write-Output "Updated $($data.Count) rows"
